$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1772.25
$ws.Range("J41").Value = 2060.4443
$ws.Range("L41").Value = 2060.4443
$ws.Range("N41").Value = -2940.4443

$ws.Range("H53").Value = 493.25
$ws.Range("I53").Value = 499.46155
$ws.Range("J53").Value = 466.33334
$ws.Range("K53").Value = 499.46155
$ws.Range("L53").Value = 466.33334
$ws.Range("M53").Value = 137.53845
$ws.Range("N53").Value = -1740.33334

$ws.Range("H55").Value = 450
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 450
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -236

$ws.Range("H76").Value = 6851.857
$ws.Range("I76").Value = 6058.6665
$ws.Range("J76").Value = 7446.75
$ws.Range("K76").Value = 6058.6665
$ws.Range("L76").Value = 7446.75
$ws.Range("M76").Value = -5743.6665
$ws.Range("N76").Value = -8076.75

$ws.Range("H79").Value = 6851.857
$ws.Range("I79").Value = 6058.6665
$ws.Range("J79").Value = 7446.75
$ws.Range("K79").Value = 6058.6665
$ws.Range("L79").Value = 7446.75
$ws.Range("M79").Value = -4966.6665
$ws.Range("N79").Value = -9630.75

$ws.Range("H100").Value = 2679.8
$ws.Range("J100").Value = 400
$ws.Range("L100").Value = 400
$ws.Range("N100").Value = -1482

$ws.Range("H135").Value = 2000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = $null
$ws.Range("M135").Value = 18000
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2444.1875
$ws.Range("I2").Value = 1735.8182
$ws.Range("J2").Value = 4002.6
$ws.Range("K2").Value = 1735.8182
$ws.Range("L2").Value = 4002.6
$ws.Range("M2").Value = -1622.8182
$ws.Range("N2").Value = -4228.6

$ws.Range("H45").Value = 1746
$ws.Range("I45").Value = 1746
$ws.Range("K45").Value = 1746
$ws.Range("M45").Value = -1369

$ws.Range("H61").Value = 4465.857
$ws.Range("I61").Value = 2106
$ws.Range("J61").Value = 5409.8
$ws.Range("K61").Value = 2106
$ws.Range("L61").Value = 5409.8
$ws.Range("M61").Value = -1894
$ws.Range("N61").Value = -5833.8

$ws.Range("H104").Value = 39994.5
$ws.Range("J104").Value = 39994.5
$ws.Range("L104").Value = 39994.5
$ws.Range("N104").Value = -46982.5

$ws.Range("H110").Value = 2979.3845
$ws.Range("I110").Value = 2045.4286
$ws.Range("J110").Value = 4069
$ws.Range("K110").Value = 2045.4286
$ws.Range("L110").Value = 4069
$ws.Range("M110").Value = -0.4285999999999603
$ws.Range("N110").Value = -8159

$ws.Range("H116").Value = 2444.1875
$ws.Range("I116").Value = 1735.8182
$ws.Range("J116").Value = 4002.6
$ws.Range("K116").Value = 1735.8182
$ws.Range("L116").Value = 4002.6
$ws.Range("M116").Value = 558.1818000000001
$ws.Range("N116").Value = -8590.6

$ws.Range("H136").Value = 4465.857
$ws.Range("I136").Value = 2106
$ws.Range("J136").Value = 5409.8
$ws.Range("K136").Value = 6318
$ws.Range("L136").Value = 16229.4
$ws.Range("M136").Value = -3768
$ws.Range("N136").Value = -21329.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2444.1875
$ws.Range("I3").Value = 1735.8182
$ws.Range("J3").Value = 4002.6
$ws.Range("K3").Value = 1735.8182
$ws.Range("L3").Value = 4002.6
$ws.Range("M3").Value = -1621.8182
$ws.Range("N3").Value = -4230.6

$ws.Range("H20").Value = 7999.25
$ws.Range("I20").Value = 8166.3335
$ws.Range("J20").Value = 7498
$ws.Range("K20").Value = 8166.3335
$ws.Range("L20").Value = 7498
$ws.Range("M20").Value = -7919.3335
$ws.Range("N20").Value = -7992

$ws.Range("H80").Value = 351
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 351
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = $null
$ws.Range("M80").Value = 351
$ws.Range("N80").Value = -2347

$ws.Range("H82").Value = 17573.555

$ws.Range("H83").Value = 351
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 351
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = $null
$ws.Range("M83").Value = 1755
$ws.Range("N83").Value = -11739

$ws.Range("H85").Value = 17573.555

$ws.Range("H86").Value = 1261.8
$ws.Range("I86").Value = 770
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 770
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = 353
$ws.Range("N86").Value = -4245.5

$ws.Range("H89").Value = 1261.8
$ws.Range("I89").Value = 770
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 3850
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = 1766
$ws.Range("N89").Value = -21229.5

$ws.Range("H132").Value = 80000
$ws.Range("J132").Value = 80000
$ws.Range("L132").Value = 80000
$ws.Range("N132").Value = -90120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6970.125
$ws.Range("I16").Value = 5624.6665
$ws.Range("J16").Value = 11006.5
$ws.Range("K16").Value = 5624.6665
$ws.Range("L16").Value = 11006.5
$ws.Range("M16").Value = -5337.6665
$ws.Range("N16").Value = -11580.5

$ws.Range("H62").Value = 9268.77
$ws.Range("I62").Value = 9942.143
$ws.Range("J62").Value = 8483.166999999999
$ws.Range("K62").Value = 9942.143
$ws.Range("L62").Value = 8483.166999999999
$ws.Range("M62").Value = -9318.143
$ws.Range("N62").Value = -9731.166999999999

$ws.Range("H65").Value = 9268.77
$ws.Range("I65").Value = 9942.143
$ws.Range("J65").Value = 8483.166999999999
$ws.Range("K65").Value = 49710.715
$ws.Range("L65").Value = 42415.835
$ws.Range("M65").Value = -46590.715
$ws.Range("N65").Value = -48655.835

$ws.Range("H99").Value = 7966.6665
$ws.Range("I99").Value = 8900
$ws.Range("K99").Value = 8900
$ws.Range("M99").Value = -7402

$ws.Range("H113").Value = 6970.125
$ws.Range("I113").Value = 5624.6665
$ws.Range("J113").Value = 11006.5
$ws.Range("K113").Value = 5624.6665
$ws.Range("L113").Value = 11006.5
$ws.Range("M113").Value = -3454.6665
$ws.Range("N113").Value = -15346.5

$ws.Range("H126").Value = 7966.6665
$ws.Range("I126").Value = 8900
$ws.Range("K126").Value = 26700
$ws.Range("M126").Value = -24230

$ws.Range("H132").Value = 1580.8
$ws.Range("I132").Value = 1538.5
$ws.Range("K132").Value = 4615.5
$ws.Range("M132").Value = -2085.5

$ws.Range("H134").Value = 4947
$ws.Range("I134").Value = 4947
$ws.Range("K134").Value = 14841
$ws.Range("M134").Value = -12306

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 6137.9287
$ws.Range("I81").Value = 6000
$ws.Range("J81").Value = 6160.9165
$ws.Range("K81").Value = 18000
$ws.Range("L81").Value = 18482.7495
$ws.Range("M81").Value = -16877
$ws.Range("N81").Value = -20728.7495

$ws.Range("H84").Value = 6137.9287
$ws.Range("I84").Value = 6000
$ws.Range("J84").Value = 6160.9165
$ws.Range("K84").Value = 54000
$ws.Range("L84").Value = 55448.2485
$ws.Range("M84").Value = -48384
$ws.Range("N84").Value = -66680.2485

$ws.Range("H86").Value = 449
$ws.Range("I86").Value = 297.5
$ws.Range("K86").Value = 892.5
$ws.Range("M86").Value = 293.5

$ws.Range("H89").Value = 449
$ws.Range("I89").Value = 297.5
$ws.Range("K89").Value = 2677.5
$ws.Range("M89").Value = 3250.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3984.5386
$ws.Range("I80").Value = 3798.5
$ws.Range("K80").Value = 3798.5
$ws.Range("M80").Value = -2800.5

$ws.Range("H83").Value = 3984.5386
$ws.Range("I83").Value = 3798.5
$ws.Range("K83").Value = 18992.5
$ws.Range("M83").Value = -14000.5

$ws.Range("H102").Value = 2861.875
$ws.Range("I102").Value = 2861.875
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2861.875
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -1239.875

$ws.Range("H122").Value = 44623.117
$ws.Range("I122").Value = 44184.645
$ws.Range("K122").Value = 132553.935
$ws.Range("M122").Value = -130103.935

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4129.8184
$ws.Range("I16").Value = 4129.8184
$ws.Range("K16").Value = 4129.8184
$ws.Range("M16").Value = -3959.8184

$ws.Range("H22").Value = 2207.1428
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 2207.1428
$ws.Range("J27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("N27").Value = -2714

$ws.Range("H40").Value = 6100
$ws.Range("I40").Value = 6100
$ws.Range("K40").Value = 6100
$ws.Range("M40").Value = -5964

$ws.Range("H46").Value = 3187
$ws.Range("J46").Value = 3274.4
$ws.Range("L46").Value = 3274.4
$ws.Range("N46").Value = -3650.4

$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

$ws.Range("H136").Value = 4188.143
$ws.Range("I136").Value = 2962.4
$ws.Range("K136").Value = 8887.200000000001
$ws.Range("M136").Value = -6337.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 6842
$ws.Range("I58").Value = 6842
$ws.Range("K58").Value = 6842
$ws.Range("M58").Value = -6534

$ws.Range("H126").Value = 4842
$ws.Range("J126").Value = 5999
$ws.Range("L126").Value = 17997
$ws.Range("N126").Value = -22937
